# Updated the quick look presentation
# Nudge "Picture 23" on slide 1 up slightly (y offset 6049934 -> 6020437 EMU)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("Picture 23")

# PowerPoint COM positions are expressed in points (1 inch = 914400 EMU = 72 pt)
$shape.Top = 6020437 / 914400 * 72
